# GUI: Stats updates JobStatus WIP
#
# The "CRUD" folder row on Sheet1 (row 2) pulls its Ready-to-Write /
# Total-Test-Suites / Total-Test-Cases counts from an external workbook
# (_Test_Suite_Statistics_for_Folders.xlsx, external reference [1]) via
#   B2 = [1]Sheet1!$H$2   (Ready to Write)
#   D2 = [1]Sheet1!$H$1   (Total Test Suites)
#   H2 = [1]Sheet1!$H$5   (Total Test Cases)
#
# The source numbers changed (0 -> 2, 33 -> 34, 256 -> 263). Refresh the
# linked values on Sheet1 so every dependent total/percentage on the
# sheet (L1, N1, E2, I2, L2, N2, N3, P3, L5, L7, ...) recalculates from
# the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 34
$ws.Range("H2").Value = 263

$excel.CalculateFull()
